$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 321, shifting existing rows 321-356 down to 322-357.
$ws.Rows(321).Insert()

# Populate the new row 321 with its data.
$ws.Range("A321").Value = 7
$ws.Range("B321").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C321").Value = "Ñuble"
$ws.Range("D321").Value = 44918
$ws.Range("E321").Value = 16
$ws.Range("F321").Value = 100114013
$ws.Range("G321").Value = "Zanahoria"
$ws.Range("H321").Value = "Sin especificar"
$ws.Range("I321").Value = "Primera"
$ws.Range("J321").Value = 160
$ws.Range("K321").Value = 8000
$ws.Range("L321").Value = 8500
$ws.Range("M321").Value = 8250
$ws.Range("N321").Value = "$/saco 20 kilos"
$ws.Range("O321").Value = "Región de Ñuble"
$ws.Range("P321").Value = 412
$ws.Range("Q321").Value = 20
$ws.Range("R321").Value = "Hortaliza"
